$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet holds a phone-number "inscription" log (Telefone / DDD / Data
# Inscrição) in rows 2..26. Three new entries were added at the top of the
# list (most recent first), pushing every existing record down by three
# rows. Insert the three rows right after the header, copying the
# formatting of the data rows below (not the header) by inserting the
# blank rows just beneath the current first three data rows, moving the
# old row 2/3/4 content down into the freshly inserted rows, and then
# writing the brand-new records into row 2/3/4.

$ws.Rows.Item(5).Insert()
$ws.Rows.Item(5).Insert()
$ws.Rows.Item(5).Insert()

# Keep these as plain text too (phone numbers / DDDs / ISO dates), so the
# relocated values below don't get reinterpreted as numbers or dates.
$ws.Range("A5:C7").NumberFormat = "@"

# Preserve the original top three records by relocating them into the
# rows that were just opened up (rows 5, 6, 7).
$ws.Range("A5").Value = $ws.Range("A2").Value2
$ws.Range("B5").Value = $ws.Range("B2").Value2
$ws.Range("C5").Value = $ws.Range("C2").Value2

$ws.Range("A6").Value = $ws.Range("A3").Value2
$ws.Range("B6").Value = $ws.Range("B3").Value2
$ws.Range("C6").Value = $ws.Range("C3").Value2

$ws.Range("A7").Value = $ws.Range("A4").Value2
$ws.Range("B7").Value = $ws.Range("B4").Value2
$ws.Range("C7").Value = $ws.Range("C4").Value2

# Force text formatting on rows 2..4 so phone numbers keep their leading
# "+" and dates stay as plain "yyyy-mm-dd" strings instead of being
# reinterpreted as numbers/dates.
$ws.Range("A2:C4").NumberFormat = "@"

# Write the three new records (newest first).
$ws.Range("A2").Value = "+5514997883211"
$ws.Range("B2").Value = "14"
$ws.Range("C2").Value = "2024-10-31"

$ws.Range("A3").Value = "+5514996538302"
$ws.Range("B3").Value = "14"
$ws.Range("C3").Value = "2024-10-25"

$ws.Range("A4").Value = "+5511964652979"
$ws.Range("B4").Value = "11"
$ws.Range("C4").Value = "2024-10-21"
